$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in zeros for every empty cell in the H3:Q20 matrix block,
# leaving the existing non-zero stoichiometry values untouched.
for ($r = 3; $r -le 20; $r++) {
    for ($c = 8; $c -le 17; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Value2 -eq $null) {
            $cell.Value2 = 0
        }
    }
}

# Update the active selection to H3:Q20 with H3 as the active cell.
$ws.Range("H3:Q20").Select()
